$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "url" column header in I1, bold to match other header cells
$ws.Range("I1").Value = 'url'
$ws.Range("I1").Font.Bold = $true

# Populate the url values for rows 2-21
$ws.Range("I2").Value = 'https://devhub.com/nibh/quisque/id/justo/sit.js?platea=morbi&dictumst=a&aliquam=ipsum&augue=integer&quam=a&sollicitudin=nibh&vitae=in&consectetuer=quis&eget=justo&rutrum=maecenas&at=rhoncus&lorem=aliquam&integer=lacus&tincidunt=morbi&ante=quis&vel=tortor&ipsum=id&praesent=nulla&blandit=ultrices&lacinia=aliquet&erat=maecenas&vestibulum=leo&sed=odio&magna=condimentum&at=id&nunc=luctus&commodo=nec&placerat=molestie&praesent=sed&blandit=justo&nam=pellentesque&nulla=viverra&integer=pede&pede=ac&justo=diam&lacinia=cras&eget=pellentesque&tincidunt=volutpat&eget=dui&tempus=maecenas&vel=tristique&pede=est&morbi=et&porttitor=tempus&lorem=semper&id=est&ligula=quam&suspendisse=pharetra&ornare=magna&consequat=ac&lectus=consequat&in=metus&est=sapien&risus=ut&auctor=nunc&sed=vestibulum&tristique=ante&in=ipsum&tempus=primis&sit=in&amet=faucibus&sem=orci&fusce=luctus&consequat=et&nulla=ultrices&nisl=posuere&nunc=cubilia&nisl=curae&duis=mauris&bibendum=viverra&felis=diam&sed=vitae&interdum=quam&venenatis=suspendisse&turpis=potenti&enim=nullam&blandit=porttitor&mi=lacus&in=at&porttitor=turpis&pede=donec&justo=posuere&eu=metus&massa=vitae&donec=ipsum&dapibus=aliquam&duis=non&at=mauris&velit=morbi&eu=non&est=lectus&congue=aliquam'
$ws.Range("I3").Value = 'https://symantec.com/ligula.xml?platea=pede&dictumst=ac&etiam=diam&faucibus=cras&cursus=pellentesque&urna=volutpat&ut=dui&tellus=maecenas&nulla=tristique&ut=est&erat=et&id=tempus&mauris=semper&vulputate=est&elementum=quam&nullam=pharetra&varius=magna&nulla=ac&facilisi=consequat&cras=metus&non=sapien&velit=ut&nec=nunc&nisi=vestibulum&vulputate=ante&nonummy=ipsum&maecenas=primis&tincidunt=in&lacus=faucibus&at=orci&velit=luctus&vivamus=et&vel=ultrices&nulla=posuere&eget=cubilia&eros=curae&elementum=mauris&pellentesque=viverra&quisque=diam&porta=vitae&volutpat=quam&erat=suspendisse&quisque=potenti&erat=nullam&eros=porttitor&viverra=lacus&eget=at&congue=turpis&eget=donec&semper=posuere&rutrum=metus&nulla=vitae&nunc=ipsum&purus=aliquam&phasellus=non&in=mauris&felis=morbi&donec=non&semper=lectus&sapien=aliquam&a=sit&libero=amet&nam=diam&dui=in&proin=magna&leo=bibendum&odio=imperdiet&porttitor=nullam&id=orci&consequat=pede&in=venenatis&consequat=non&ut=sodales&nulla=sed&sed=tincidunt&accumsan=eu&felis=felis&ut=fusce&at=posuere'
$ws.Range("I4").Value = 'https://samsung.com/lorem/id/ligula/suspendisse/ornare.json?aenean=sed&auctor=magna&gravida=at&sem=nunc&praesent=commodo&id=placerat&massa=praesent&id=blandit&nisl=nam&venenatis=nulla&lacinia=integer&aenean=pede&sit=justo&amet=lacinia&justo=eget&morbi=tincidunt&ut=eget&odio=tempus&cras=vel&mi=pede&pede=morbi&malesuada=porttitor&in=lorem&imperdiet=id&et=ligula&commodo=suspendisse&vulputate=ornare&justo=consequat&in=lectus&blandit=in&ultrices=est&enim=risus&lorem=auctor&ipsum=sed&dolor=tristique&sit=in&amet=tempus&consectetuer=sit&adipiscing=amet&elit=sem&proin=fusce&interdum=consequat&mauris=nulla&non=nisl&ligula=nunc&pellentesque=nisl&ultrices=duis&phasellus=bibendum&id=felis&sapien=sed&in=interdum&sapien=venenatis&iaculis=turpis&congue=enim&vivamus=blandit&metus=mi&arcu=in&adipiscing=porttitor&molestie=pede&hendrerit=justo&at=eu&vulputate=massa&vitae=donec&nisl=dapibus&aenean=duis&lectus=at&pellentesque=velit&eget=eu&nunc=est&donec=congue&quis=elementum&orci=in&eget=hac&orci=habitasse&vehicula=platea&condimentum=dictumst&curabitur=morbi&in=vestibulum&libero=velit&ut=id&massa=pretium&volutpat=iaculis&convallis=diam&morbi=erat&odio=fermentum&odio=justo&elementum=nec&eu=condimentum&interdum=neque&eu=sapien&tincidunt=placerat&in=ante&leo=nulla&maecenas=justo&pulvinar=aliquam&lobortis=quis'
$ws.Range("I5").Value = 'https://youtu.be/accumsan/tortor/quis/turpis/sed.js?risus=bibendum&auctor=imperdiet&sed=nullam&tristique=orci&in=pede&tempus=venenatis&sit=non&amet=sodales&sem=sed&fusce=tincidunt&consequat=eu&nulla=felis&nisl=fusce&nunc=posuere&nisl=felis&duis=sed&bibendum=lacus&felis=morbi&sed=sem'
$ws.Range("I6").Value = 'https://cbsnews.com/elit/ac.xml?magnis=faucibus&dis=orci&parturient=luctus&montes=et&nascetur=ultrices&ridiculus=posuere&mus=cubilia&vivamus=curae&vestibulum=donec&sagittis=pharetra&sapien=magna&cum=vestibulum&sociis=aliquet&natoque=ultrices&penatibus=erat&et=tortor&magnis=sollicitudin&dis=mi&parturient=sit&montes=amet&nascetur=lobortis&ridiculus=sapien&mus=sapien&etiam=non&vel=mi&augue=integer&vestibulum=ac&rutrum=neque&rutrum=duis&neque=bibendum&aenean=morbi&auctor=non&gravida=quam&sem=nec'
$ws.Range("I7").Value = 'https://zimbio.com/et/tempus/semper/est/quam/pharetra/magna.jsp?sit=sit&amet=amet&eros=lobortis&suspendisse=sapien&accumsan=sapien&tortor=non&quis=mi&turpis=integer&sed=ac&ante=neque&vivamus=duis&tortor=bibendum&duis=morbi&mattis=non&egestas=quam&metus=nec&aenean=dui&fermentum=luctus&donec=rutrum&ut=nulla&mauris=tellus&eget=in&massa=sagittis&tempor=dui&convallis=vel&nulla=nisl&neque=duis&libero=ac&convallis=nibh&eget=fusce&eleifend=lacus&luctus=purus&ultricies=aliquet&eu=at&nibh=feugiat&quisque=non&id=pretium&justo=quis&sit=lectus&amet=suspendisse&sapien=potenti&dignissim=in&vestibulum=eleifend&vestibulum=quam&ante=a&ipsum=odio&primis=in&in=hac&faucibus=habitasse&orci=platea&luctus=dictumst&et=maecenas&ultrices=ut&posuere=massa&cubilia=quis&curae=augue&nulla=luctus&dapibus=tincidunt&dolor=nulla&vel=mollis&est=molestie&donec=lorem&odio=quisque&justo=ut&sollicitudin=erat&ut=curabitur&suscipit=gravida&a=nisi&feugiat=at&et=nibh&eros=in&vestibulum=hac&ac=habitasse&est=platea&lacinia=dictumst&nisi=aliquam&venenatis=augue&tristique=quam&fusce=sollicitudin&congue=vitae&diam=consectetuer&id=eget&ornare=rutrum&imperdiet=at&sapien=lorem&urna=integer&pretium=tincidunt&nisl=ante&ut=vel&volutpat=ipsum&sapien=praesent&arcu=blandit&sed=lacinia'
$ws.Range("I8").Value = 'https://bbc.co.uk/pulvinar/lobortis/est/phasellus/sit/amet/erat.aspx?phasellus=ligula&in=in&felis=lacus&donec=curabitur&semper=at&sapien=ipsum&a=ac&libero=tellus&nam=semper&dui=interdum&proin=mauris&leo=ullamcorper&odio=purus&porttitor=sit&id=amet&consequat=nulla&in=quisque&consequat=arcu&ut=libero&nulla=rutrum'
$ws.Range("I9").Value = 'http://hostgator.com/at/dolor/quis/odio/consequat/varius.jsp?ultricies=integer&eu=tincidunt&nibh=ante&quisque=vel&id=ipsum&justo=praesent&sit=blandit&amet=lacinia&sapien=erat'
$ws.Range("I10").Value = 'http://whitehouse.gov/sapien/varius/ut/blandit.aspx?dis=lectus&parturient=in&montes=quam&nascetur=fringilla&ridiculus=rhoncus&mus=mauris&vivamus=enim&vestibulum=leo&sagittis=rhoncus&sapien=sed&cum=vestibulum&sociis=sit&natoque=amet&penatibus=cursus&et=id&magnis=turpis&dis=integer&parturient=aliquet&montes=massa&nascetur=id&ridiculus=lobortis&mus=convallis&etiam=tortor&vel=risus&augue=dapibus&vestibulum=augue&rutrum=vel&rutrum=accumsan&neque=tellus&aenean=nisi&auctor=eu&gravida=orci&sem=mauris&praesent=lacinia&id=sapien&massa=quis&id=libero&nisl=nullam&venenatis=sit&lacinia=amet&aenean=turpis&sit=elementum&amet=ligula&justo=vehicula&morbi=consequat&ut=morbi&odio=a&cras=ipsum&mi=integer&pede=a&malesuada=nibh&in=in&imperdiet=quis&et=justo&commodo=maecenas&vulputate=rhoncus&justo=aliquam&in=lacus&blandit=morbi&ultrices=quis&enim=tortor&lorem=id&ipsum=nulla&dolor=ultrices&sit=aliquet&amet=maecenas&consectetuer=leo&adipiscing=odio&elit=condimentum&proin=id&interdum=luctus&mauris=nec&non=molestie&ligula=sed&pellentesque=justo&ultrices=pellentesque&phasellus=viverra&id=pede&sapien=ac&in=diam&sapien=cras&iaculis=pellentesque&congue=volutpat&vivamus=dui&metus=maecenas&arcu=tristique&adipiscing=est&molestie=et&hendrerit=tempus&at=semper&vulputate=est'
$ws.Range("I11").Value = 'https://sciencedirect.com/fermentum/justo/nec/condimentum/neque.jpg?duis=elementum&faucibus=ligula&accumsan=vehicula&odio=consequat&curabitur=morbi&convallis=a&duis=ipsum&consequat=integer&dui=a&nec=nibh&nisi=in&volutpat=quis&eleifend=justo&donec=maecenas&ut=rhoncus&dolor=aliquam&morbi=lacus&vel=morbi&lectus=quis&in=tortor&quam=id&fringilla=nulla&rhoncus=ultrices&mauris=aliquet&enim=maecenas&leo=leo&rhoncus=odio&sed=condimentum&vestibulum=id&sit=luctus&amet=nec&cursus=molestie&id=sed&turpis=justo&integer=pellentesque&aliquet=viverra&massa=pede&id=ac&lobortis=diam&convallis=cras&tortor=pellentesque&risus=volutpat&dapibus=dui&augue=maecenas&vel=tristique&accumsan=est&tellus=et&nisi=tempus&eu=semper&orci=est&mauris=quam&lacinia=pharetra&sapien=magna&quis=ac&libero=consequat&nullam=metus&sit=sapien&amet=ut&turpis=nunc&elementum=vestibulum&ligula=ante&vehicula=ipsum&consequat=primis&morbi=in&a=faucibus&ipsum=orci&integer=luctus&a=et&nibh=ultrices&in=posuere'
$ws.Range("I12").Value = 'https://home.pl/ipsum/primis/in.xml?volutpat=morbi'
$ws.Range("I13").Value = 'https://hhs.gov/sodales/scelerisque/mauris.js?erat=duis&tortor=aliquam&sollicitudin=convallis&mi=nunc&sit=proin&amet=at&lobortis=turpis&sapien=a&sapien=pede&non=posuere&mi=nonummy&integer=integer&ac=non&neque=velit&duis=donec&bibendum=diam&morbi=neque&non=vestibulum&quam=eget&nec=vulputate&dui=ut&luctus=ultrices&rutrum=vel&nulla=augue&tellus=vestibulum&in=ante&sagittis=ipsum&dui=primis&vel=in&nisl=faucibus&duis=orci&ac=luctus'
$ws.Range("I14").Value = 'https://cloudflare.com/imperdiet/et/commodo/vulputate.xml?sit=sed&amet=tincidunt&eros=eu&suspendisse=felis&accumsan=fusce&tortor=posuere&quis=felis&turpis=sed&sed=lacus&ante=morbi&vivamus=sem&tortor=mauris&duis=laoreet&mattis=ut&egestas=rhoncus&metus=aliquet&aenean=pulvinar&fermentum=sed&donec=nisl&ut=nunc&mauris=rhoncus&eget=dui&massa=vel&tempor=sem&convallis=sed&nulla=sagittis&neque=nam&libero=congue&convallis=risus&eget=semper&eleifend=porta&luctus=volutpat&ultricies=quam&eu=pede&nibh=lobortis&quisque=ligula&id=sit&justo=amet&sit=eleifend&amet=pede&sapien=libero&dignissim=quis&vestibulum=orci&vestibulum=nullam&ante=molestie'
$ws.Range("I15").Value = 'http://lycos.com/erat/id/mauris/vulputate.xml?iaculis=ullamcorper&congue=augue&vivamus=a&metus=suscipit&arcu=nulla&adipiscing=elit&molestie=ac'
$ws.Range("I16").Value = 'http://twitpic.com/nisi.json?lacinia=non&eget=sodales&tincidunt=sed&eget=tincidunt&tempus=eu&vel=felis&pede=fusce&morbi=posuere&porttitor=felis&lorem=sed&id=lacus&ligula=morbi&suspendisse=sem&ornare=mauris&consequat=laoreet&lectus=ut&in=rhoncus&est=aliquet&risus=pulvinar&auctor=sed&sed=nisl&tristique=nunc&in=rhoncus&tempus=dui&sit=vel&amet=sem&sem=sed&fusce=sagittis&consequat=nam&nulla=congue&nisl=risus&nunc=semper&nisl=porta&duis=volutpat&bibendum=quam&felis=pede&sed=lobortis&interdum=ligula&venenatis=sit&turpis=amet&enim=eleifend&blandit=pede&mi=libero&in=quis&porttitor=orci&pede=nullam&justo=molestie&eu=nibh&massa=in&donec=lectus&dapibus=pellentesque&duis=at&at=nulla&velit=suspendisse&eu=potenti&est=cras&congue=in&elementum=purus&in=eu&hac=magna&habitasse=vulputate&platea=luctus&dictumst=cum&morbi=sociis&vestibulum=natoque&velit=penatibus&id=et&pretium=magnis'
$ws.Range("I17").Value = 'https://nsw.gov.au/vestibulum/ac/est/lacinia.js?nunc=ante'
$ws.Range("I18").Value = 'http://imdb.com/vel/est/donec/odio/justo.aspx?tortor=nulla&quis=neque&turpis=libero&sed=convallis&ante=eget&vivamus=eleifend&tortor=luctus&duis=ultricies&mattis=eu&egestas=nibh&metus=quisque&aenean=id&fermentum=justo&donec=sit&ut=amet&mauris=sapien&eget=dignissim&massa=vestibulum&tempor=vestibulum&convallis=ante&nulla=ipsum&neque=primis&libero=in&convallis=faucibus&eget=orci&eleifend=luctus&luctus=et&ultricies=ultrices&eu=posuere&nibh=cubilia&quisque=curae&id=nulla&justo=dapibus'
$ws.Range("I19").Value = 'https://smugmug.com/in/tempus.jpg?sit=faucibus&amet=cursus&consectetuer=urna&adipiscing=ut&elit=tellus&proin=nulla&risus=ut&praesent=erat&lectus=id&vestibulum=mauris&quam=vulputate&sapien=elementum&varius=nullam&ut=varius&blandit=nulla&non=facilisi&interdum=cras&in=non&ante=velit&vestibulum=nec&ante=nisi&ipsum=vulputate&primis=nonummy&in=maecenas&faucibus=tincidunt&orci=lacus&luctus=at&et=velit&ultrices=vivamus&posuere=vel&cubilia=nulla&curae=eget&duis=eros&faucibus=elementum&accumsan=pellentesque&odio=quisque&curabitur=porta&convallis=volutpat&duis=erat&consequat=quisque&dui=erat&nec=eros&nisi=viverra&volutpat=eget&eleifend=congue&donec=eget&ut=semper&dolor=rutrum&morbi=nulla&vel=nunc&lectus=purus&in=phasellus&quam=in&fringilla=felis&rhoncus=donec&mauris=semper&enim=sapien&leo=a&rhoncus=libero&sed=nam&vestibulum=dui&sit=proin&amet=leo&cursus=odio&id=porttitor&turpis=id&integer=consequat&aliquet=in&massa=consequat&id=ut&lobortis=nulla&convallis=sed&tortor=accumsan&risus=felis'
$ws.Range("I20").Value = 'http://simplemachines.org/erat/curabitur/gravida.jpg?convallis=neque&nunc=libero&proin=convallis&at=eget&turpis=eleifend&a=luctus&pede=ultricies&posuere=eu&nonummy=nibh&integer=quisque&non=id&velit=justo&donec=sit&diam=amet&neque=sapien&vestibulum=dignissim&eget=vestibulum&vulputate=vestibulum&ut=ante&ultrices=ipsum&vel=primis&augue=in&vestibulum=faucibus&ante=orci&ipsum=luctus&primis=et&in=ultrices&faucibus=posuere&orci=cubilia&luctus=curae&et=nulla&ultrices=dapibus&posuere=dolor&cubilia=vel&curae=est&donec=donec&pharetra=odio&magna=justo&vestibulum=sollicitudin&aliquet=ut&ultrices=suscipit&erat=a&tortor=feugiat&sollicitudin=et&mi=eros&sit=vestibulum&amet=ac&lobortis=est&sapien=lacinia&sapien=nisi&non=venenatis&mi=tristique&integer=fusce&ac=congue&neque=diam&duis=id'
$ws.Range("I21").Value = 'https://spiegel.de/vitae/ipsum/aliquam/non/mauris/morbi/non.png?hendrerit=convallis&at=tortor&vulputate=risus&vitae=dapibus&nisl=augue&aenean=vel&lectus=accumsan&pellentesque=tellus&eget=nisi&nunc=eu&donec=orci&quis=mauris&orci=lacinia&eget=sapien&orci=quis&vehicula=libero&condimentum=nullam&curabitur=sit&in=amet&libero=turpis&ut=elementum&massa=ligula&volutpat=vehicula&convallis=consequat&morbi=morbi&odio=a&odio=ipsum&elementum=integer&eu=a&interdum=nibh&eu=in&tincidunt=quis&in=justo&leo=maecenas&maecenas=rhoncus&pulvinar=aliquam&lobortis=lacus&est=morbi&phasellus=quis&sit=tortor&amet=id&erat=nulla&nulla=ultrices&tempus=aliquet&vivamus=maecenas&in=leo&felis=odio&eu=condimentum&sapien=id&cursus=luctus&vestibulum=nec&proin=molestie&eu=sed&mi=justo&nulla=pellentesque&ac=viverra&enim=pede&in=ac&tempor=diam&turpis=cras&nec=pellentesque&euismod=volutpat&scelerisque=dui&quam=maecenas&turpis=tristique&adipiscing=est&lorem=et&vitae=tempus&mattis=semper&nibh=est&ligula=quam&nec=pharetra&sem=magna&duis=ac&aliquam=consequat&convallis=metus&nunc=sapien&proin=ut&at=nunc&turpis=vestibulum&a=ante&pede=ipsum&posuere=primis&nonummy=in&integer=faucibus&non=orci&velit=luctus&donec=et'

# Update the active selection to match the author's final selection
$ws.Range("I2:I21").Select() | Out-Null
